$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2, column A: registration number becomes bold (no centering, no border) ---
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Value = "F22017001023"

# --- Row 2, column B: new semester number, keep it as centered quoted text ---
$ws.Range("B2").Value = "'5"

# --- Row 1, columns C & D: new headers using Consolas 10pt, text number format ---
$ws.Range("C1").Value = "Subject"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Font.Size = 10
$ws.Range("C1").Font.Name = "Consolas"

$ws.Range("D1").Value = "Sem Mark"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Font.Size = 10
$ws.Range("D1").Font.Name = "Consolas"

# --- Row 2, column C: subject name, bordered & centered, with a "no-fill" fill applied ---
$ws.Range("C2").Value = "TH5"
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").Interior.ColorIndex = -4142

# --- Row 2, column D: numeric mark, bordered & centered ---
$ws.Range("D2").Value = 22
$ws.Range("D2").Borders.LineStyle = 1
$ws.Range("D2").HorizontalAlignment = -4108

# --- Column width for the new "Subject" column ---
$ws.Columns.Item(3).ColumnWidth = 34.5

# --- Selection follows the newly entered mark cell ---
$ws.Range("D2").Select()
